# Updating coverage values for sch
# The "coverage" values for every even year (2018, 2020, 2022, ...) on the
# "Platform Coverage" sheet move from 0.75 down to 0.6.
$wb = $excel.ActiveWorkbook
$wsCoverage = $wb.Worksheets.Item("Platform Coverage")

$cols = @("H","J","L","N","P","R","T","V","X","Z","AB","AD")
foreach ($col in $cols) {
    $wsCoverage.Range($col + "2").Value = 0.6
}

# Bring the "Platform Coverage" sheet to the front and leave the selection
# parked on P9, scrolled so column M is the left-most visible column -
# mirroring where the author was working when they made the edit.
$wsCoverage.Activate()
$wsCoverage.Application.ActiveWindow.ScrollColumn = 13
$wsCoverage.Range("P9").Select()

$wb.Save()
